$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: UserID / Email
$ws.Cells.Item(2, 1).Value = 55656862326
$ws.Cells.Item(2, 2).Value = "HJVKJ@GMAIL.COM"

# Widen column A so the UserID values are fully visible
$ws.Columns.Item(1).ColumnWidth = 14.65

# Leave the selection where the user ended up after entering the data
$ws.Range("A4").Select()
